$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "36.749.59"
$ws.Range("E2").Value = "  +3.86%  "
$ws.Range("D3").Value = "1.909.11"
$ws.Range("E3").Value = "  +1.36%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Formula = "=""249.16"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  +1.06%  "
$ws.Range("D6").Formula = "=""0.696"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Formula = "=""46.64"""
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "  +7.25%  "
$ws.Range("E9").Value = "  +4.91%  "
$ws.Range("D10").Formula = "=""57.90"""
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  +7.90%  "
$ws.Range("D11").Formula = "=""0.0755"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("E12").Value = "  +2.23%  "
$ws.Range("D13").Formula = "=""14.64"""
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  +8.30%  "
$ws.Range("D14").Formula = "=""0.811"""
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  +4.81%  "
$ws.Range("D15").Value = "2.187.17"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("E16").Value = "  +2.71%  "
$ws.Range("D17").Value = "1.917.07"
$ws.Range("E17").Value = "  +2.01%  "
$ws.Range("D18").Value = "36.731.38"
$ws.Range("E18").Value = "  +3.81%  "
$ws.Range("D19").Formula = "=""74.25"""
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("E20").Value = "  +2.90%  "
$ws.Range("E21").Value = "  +5.56%  "
$ws.Range("D22").Formula = "=""250.47"""
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  +2.30%  "
$ws.Range("D23").Formula = "=""5.12"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  -1.11%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").Formula = "=""2.49"""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  -4.43%  "
$ws.Range("E26").Value = "  +1.75%  "
$ws.Range("D27").Formula = "=""166.78"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  +1.28%  "
$ws.Range("D28").Formula = "=""8.74"""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  +0.94%  "
$ws.Range("D29").Formula = "=""18.63"""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("D30").Formula = "=""0.128"""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("D31").Formula = "=""4.62"""
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  +7.76%  "
$ws.Range("D32").Formula = "=""0.0613"""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  +2.84%  "
$ws.Range("E33").Value = "  +2.84%  "
$ws.Range("E34").Value = "  +1.46%  "
$ws.Range("B35").Value = "BinanceUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D35").Formula = "=""1.00"""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Formula = "=""0.0878"""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  +19.50%  "
$ws.Range("D37").Formula = "=""18.81"""
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  +55.25%  "
$ws.Range("E38").Value = "  -0.69%  "
$ws.Range("D39").Formula = "=""0.869"""
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  +1.88%  "
$ws.Range("E40").Value = "  +1.33%  "
$ws.Range("D41").Formula = "=""104.37"""
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  +7.13%  "
$ws.Range("E42").Value = "  +4.25%  "
$ws.Range("D43").Formula = "=""17.79"""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  +2.81%  "
$ws.Range("D44").Formula = "=""2.81"""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  +17.84%  "
$ws.Range("E45").Value = "  +1.12%  "
$ws.Range("D46").Value = "1.343.90"
$ws.Range("E46").Value = "  +2.62%  "
$ws.Range("D47").Formula = "=""2.37"""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  -1.12%  "
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("E49").Value = "  +3.10%  "
$ws.Range("E50").Value = "  +1.68%  "
$ws.Range("D51").Value = "2.083.39"
$ws.Range("E51").Value = "  +0.99%  "
